$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price values are plain numeric-looking strings need to be
# forced to Text format first, otherwise Excel auto-converts them to numbers
# (losing the original formatted string / introducing float rounding).
$numericPriceCells = @("D5","D6","D10","D12","D17","D19","D20","D22","D23","D26","D28","D33","D36","D37","D38","D40","D44","D45","D47","D48","D49","D50","D51")
foreach ($addr in $numericPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "61.223.98"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "3.369.11"
$ws.Range("E3").Value = "  +1.58%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "571.18"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").Value = "137.25"
$ws.Range("E6").Value = "  +7.89%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.368.60"
$ws.Range("E8").Value = "  +1.53%  "
$ws.Range("E9").Value = "  -0.47%  "
$ws.Range("D10").Value = "7.59"
$ws.Range("E10").Value = "  +5.57%  "
$ws.Range("E11").Value = "  +4.26%  "
$ws.Range("D12").Value = "0.391"
$ws.Range("E12").Value = "  +4.34%  "
$ws.Range("D13").Value = "3.949.20"
$ws.Range("E13").Value = "  +1.76%  "
$ws.Range("E14").Value = "  +2.33%  "
$ws.Range("D15").Value = "3.376.01"
$ws.Range("E15").Value = "  +1.77%  "
$ws.Range("E16").Value = "  +1.93%  "
$ws.Range("D17").Value = "25.23"
$ws.Range("E17").Value = "  +1.62%  "
$ws.Range("D18").Value = "61.327.99"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "13.98"
$ws.Range("E19").Value = "  +6.17%  "
$ws.Range("D20").Value = "5.81"
$ws.Range("E20").Value = "  +4.09%  "
$ws.Range("E21").Value = "  +3.57%  "
$ws.Range("D22").Value = "380.46"
$ws.Range("E22").Value = "  +7.05%  "
$ws.Range("D23").Value = "0.575"
$ws.Range("E23").Value = "  +3.76%  "
$ws.Range("D24").Value = "3.507.97"
$ws.Range("E24").Value = "  +1.79%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "70.85"
$ws.Range("E26").Value = "  +0.91%  "
$ws.Range("E27").Value = "  +10.62%  "
$ws.Range("D28").Value = "1.65"
$ws.Range("E28").Value = "  +11.94%  "
$ws.Range("E29").Value = "  +7.77%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E32").Value = "  +4.44%  "
$ws.Range("D33").Value = "2.13"
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("D35").Value = "3.403.49"
$ws.Range("E35").Value = "  +1.87%  "
$ws.Range("D36").Value = "23.45"
$ws.Range("E36").Value = "  +5.02%  "
$ws.Range("D37").Value = "5.56"
$ws.Range("E37").Value = "  +1.55%  "
$ws.Range("D38").Value = "6.95"
$ws.Range("E38").Value = "  +3.27%  "
$ws.Range("E39").Value = "  +3.49%  "
$ws.Range("D40").Value = "163.70"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("E41").Value = "  +5.74%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("E43").Value = "  +3.99%  "
$ws.Range("D44").Value = "41.50"
$ws.Range("E44").Value = "  +1.67%  "
$ws.Range("D45").Value = "0.760"
$ws.Range("E45").Value = "  +1.62%  "
$ws.Range("E46").Value = "  +7.56%  "
$ws.Range("D47").Value = "1.63"
$ws.Range("E47").Value = "  +5.66%  "
$ws.Range("D48").Value = "23.29"
$ws.Range("E48").Value = "  +3.98%  "
$ws.Range("D49").Value = "6.96"
$ws.Range("E49").Value = "  +5.34%  "
$ws.Range("D50").Value = "23.25"
$ws.Range("E50").Value = "  +11.73%  "
$ws.Range("D51").Value = "2.42"
$ws.Range("E51").Value = "  +11.15%  "
